$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="25.962.22"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '="  +0.41%  "'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("D3").Formula = '="1.642.83"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = '="  +0.31%  "'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

$ws.Range("D4").Formula = '="1.002"'
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Formula = '="  +0.18%  "'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("D5").Formula = '="215.32"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = '="  +0.20%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

$ws.Range("D6").Formula = '="0.5074"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = '="  +1.15%  "'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range("D7").Formula = '="1.005"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Formula = '="  +0.47%  "'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)

$ws.Range("D8").Formula = '="0.2562"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Formula = '="  -0.20%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)

$ws.Range("D9").Formula = '="0.06374"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="  +0.12%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)

$ws.Range("D10").Formula = '="19.53"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Formula = '="  +0.03%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

$ws.Range("D11").Formula = '="0.07776"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = '="  +0.71%  "'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

$ws.Range("D12").Formula = '="4.291"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="  +1.05%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

$ws.Range("D13").Formula = '="1.648.66"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="  +0.67%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

$ws.Range("D14").Formula = '="0.5461"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Formula = '="  +0.30%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

$ws.Range("D15").Formula = '="0.0₅7834"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = '="  -0.57%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)

$ws.Range("D16").Formula = '="64.42"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = '="  +0.79%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

$ws.Range("D17").Formula = '="26.003.03"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Formula = '="  +0.53%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)

$ws.Range("E18").Formula = '="  +0.43%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)

$ws.Range("D19").Formula = '="197.62"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = '="  -2.33%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

$ws.Range("E20").Formula = '="  +1.46%  "'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)

$ws.Range("D21").Formula = '="9.972"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Formula = '="  +1.18%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

$ws.Range("D22").Formula = '="6.054"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Formula = '="  +1.48%  "'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)

$ws.Range("D23").Formula = '="1.006"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Formula = '="  +0.51%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)

$ws.Range("D24").Formula = '="1.897"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Formula = '="  +2.20%  "'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

$ws.Range("D25").Formula = '="141.49"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = '="  +0.74%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)

$ws.Range("D26").Formula = '="0.1174"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Formula = '="  +3.77%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)

$ws.Range("D27").Formula = '="6.885"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = '="  +1.81%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

$ws.Range("E28").Formula = '="  +0.31%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)

$ws.Range("D29").Formula = '="1.239"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Formula = '="  +0.23%  "'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)

$ws.Range("D30").Formula = '="0.04992"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Formula = '="  +0.77%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)

$ws.Range("D31").Formula = '="3.264"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Formula = '="  +0.06%  "'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)

$ws.Range("E32").Formula = '="  +0.18%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)

$ws.Range("D33").Formula = '="1.543"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Formula = '="  +0.34%  "'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)

$ws.Range("D34").Formula = '="2.362"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Formula = '="  +0.15%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)

$ws.Range("D35").Formula = '="0.8952"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Formula = '="  +0.75%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)

$ws.Range("D36").Formula = '="2.585"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Formula = '="  -1.41%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)

$ws.Range("D37").Formula = '="1.132.29"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = '="  -0.87%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)

$ws.Range("D38").Formula = '="0.5453"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Formula = '="  -2.99%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)

$ws.Range("D39").Formula = '="0.01554"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Formula = '="  -0.22%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)

$ws.Range("D40").Formula = '="2.557"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = '="  -0.55%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)

$ws.Range("D41").Formula = '="1.005"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Formula = '="  +0.61%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)

$ws.Range("D42").Formula = '="0.0₈128"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = '="  +10.57%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)

$ws.Range("D43").Formula = '="5.606"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = '="  -0.93%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)

$ws.Range("D44").Formula = '="0.8184"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = '="  +1.92%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)

$ws.Range("D45").Formula = '="99.84"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = '="  +0.14%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)

$ws.Range("D46").Formula = '="1.777.94"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Formula = '="  +0.13%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)

$ws.Range("D47").Formula = '="0.4547"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = '="  +0.50%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

$ws.Range("D48").Formula = '="1.004"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = '="  -0.01%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)

$ws.Range("D49").Formula = '="54.87"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Formula = '="  +0.52%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)

$ws.Range("D50").Formula = '="0.05073"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Formula = '="  +0.55%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)

$ws.Range("D51").Formula = '="1.004"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Formula = '="  +0.38%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)

$excel.CutCopyMode = $false
